# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" worksheets to the newly scraped values.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row -> new value for column F
$updates = @{
    2  = 666
    3  = 506
    6  = 49
    8  = 2637
    9  = 4158
    10 = 101
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
